# Applies the cryptos-list price/volume refresh described by the commit.
# Cells are plain inline/shared strings (General format) in the source sheet,
# so any replacement value that LOOKS like a plain number (e.g. "0.724")
# must be forced to stay text: we briefly mark the cell as Text (@), assign
# the string, then restore the "Normal" style so no stray number format
# sticks around on the cell (matches the original, style-less cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.350.39"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.914.47"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.724"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.356"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0756"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0992"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "2.194.11"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.719"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.25%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.923.50"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "35.385.70"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("D20").Value = "0.0₃0844"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.76%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("E30").Value = "  +5.66%  "
$ws.Range("D31").Value = "4.128.51"
$ws.Range("E31").Value = "  +19.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +23.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0581"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.925"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0219"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.42%  "
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").Value = "1.340.70"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.16%  "
